{"js": "// Update the two-digit division worksheet: replace each \"A\u00f7B=\" expression\n// in the table with its new value, per the authoring diff.\n// Pairs are [originalText, newText] and are matched exactly (matchCase)\n// against the unique \"A\u00f7B=\" cell contents already in the document, so\n// processing them in this fixed order is safe even though one new value\n// (\"61\u00f73=\") happens to equal a *different* cell's original value.\nconst replacements = [\n  [\"58\u00f73=\", \"71\u00f76=\"],\n  [\"96\u00f73=\", \"18\u00f73=\"],\n  [\"92\u00f73=\", \"34\u00f79=\"],\n  [\"58\u00f74=\", \"30\u00f73=\"],\n  [\"21\u00f75=\", \"72\u00f73=\"],\n  [\"10\u00f72=\", \"59\u00f78=\"],\n  [\"71\u00f72=\", \"50\u00f79=\"],\n  [\"64\u00f72=\", \"41\u00f74=\"],\n  [\"39\u00f72=\", \"16\u00f76=\"],\n  [\"73\u00f74=\", \"10\u00f74=\"],\n  [\"80\u00f73=\", \"48\u00f77=\"],\n  [\"57\u00f72=\", \"52\u00f72=\"],\n  [\"37\u00f76=\", \"87\u00f72=\"],\n  [\"61\u00f73=\", \"15\u00f72=\"],\n  [\"93\u00f75=\", \"75\u00f78=\"],\n  [\"75\u00f73=\", \"65\u00f75=\"],\n  [\"98\u00f75=\", \"30\u00f72=\"],\n  [\"67\u00f79=\", \"61\u00f73=\"],\n  [\"18\u00f77=\", \"81\u00f74=\"],\n  [\"77\u00f75=\", \"78\u00f78=\"],\n  [\"50\u00f76=\", \"60\u00f73=\"],\n  [\"49\u00f76=\", \"24\u00f77=\"],\n  [\"51\u00f72=\", \"50\u00f73=\"],\n  [\"40\u00f78=\", \"34\u00f73=\"],\n  [\"13\u00f73=\", \"70\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit division worksheet: replace each \"A\u00f7B=\" expression\n# in the table with its new value, per the authoring diff.\n# Each (old, new) pair is matched exactly against the unique \"A\u00f7B=\" text\n# already present in the document, so applying them in this fixed order is\n# safe even though one new value (\"61\u00f73=\") happens to equal a *different*\n# cell's original value.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"58\u00f73=\", \"71\u00f76=\"),\n  @(\"96\u00f73=\", \"18\u00f73=\"),\n  @(\"92\u00f73=\", \"34\u00f79=\"),\n  @(\"58\u00f74=\", \"30\u00f73=\"),\n  @(\"21\u00f75=\", \"72\u00f73=\"),\n  @(\"10\u00f72=\", \"59\u00f78=\"),\n  @(\"71\u00f72=\", \"50\u00f79=\"),\n  @(\"64\u00f72=\", \"41\u00f74=\"),\n  @(\"39\u00f72=\", \"16\u00f76=\"),\n  @(\"73\u00f74=\", \"10\u00f74=\"),\n  @(\"80\u00f73=\", \"48\u00f77=\"),\n  @(\"57\u00f72=\", \"52\u00f72=\"),\n  @(\"37\u00f76=\", \"87\u00f72=\"),\n  @(\"61\u00f73=\", \"15\u00f72=\"),\n  @(\"93\u00f75=\", \"75\u00f78=\"),\n  @(\"75\u00f73=\", \"65\u00f75=\"),\n  @(\"98\u00f75=\", \"30\u00f72=\"),\n  @(\"67\u00f79=\", \"61\u00f73=\"),\n  @(\"18\u00f77=\", \"81\u00f74=\"),\n  @(\"77\u00f75=\", \"78\u00f78=\"),\n  @(\"50\u00f76=\", \"60\u00f73=\"),\n  @(\"49\u00f76=\", \"24\u00f77=\"),\n  @(\"51\u00f72=\", \"50\u00f73=\"),\n  @(\"40\u00f78=\", \"34\u00f73=\"),\n  @(\"13\u00f73=\", \"70\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
